# Update the "想去人数" (want-to-go count) figures in column F on the
# "展览" (Exhibition) and "全部类型" (All types) sheets to match the
# newly scraped output.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 -----------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("F4").Value = 559
$ws1.Range("F5").Value = 8612
$ws1.Range("F9").Value = 6082
$ws1.Range("F13").Value = 8643
$ws1.Range("F14").Value = 10233
$ws1.Range("F15").Value = 1183
$ws1.Range("F16").Value = 1028
$ws1.Range("F17").Value = 4740
$ws1.Range("F19").Value = 369
$ws1.Range("F21").Value = 310
$ws1.Range("F22").Value = 171
$ws1.Range("F25").Value = 1837
$ws1.Range("F27").Value = 1118
$ws1.Range("F28").Value = 833
$ws1.Range("F29").Value = 1968
$ws1.Range("F31").Value = 548
$ws1.Range("F32").Value = 2504
$ws1.Range("F34").Value = 147
$ws1.Range("F41").Value = 3163
$ws1.Range("F43").Value = 74
$ws1.Range("F45").Value = 553
$ws1.Range("F48").Value = 214

# --- Sheet 4: 全部类型 ---------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)

$ws4.Range("F4").Value = 559
$ws4.Range("F5").Value = 8612
$ws4.Range("F9").Value = 6082
$ws4.Range("F11").Value = 8643
$ws4.Range("F12").Value = 10233
$ws4.Range("F14").Value = 1184
$ws4.Range("F15").Value = 1028
$ws4.Range("F16").Value = 4740
$ws4.Range("F18").Value = 369
$ws4.Range("F20").Value = 310
$ws4.Range("F22").Value = 171
$ws4.Range("F25").Value = 1837
$ws4.Range("F27").Value = 833
$ws4.Range("F29").Value = 1968
$ws4.Range("F31").Value = 548
$ws4.Range("F32").Value = 2504
$ws4.Range("F43").Value = 74
$ws4.Range("F45").Value = 553
$ws4.Range("F47").Value = 214
